$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.288.24"
$ws.Range("E2").Value = "  +2.13%  "

# Row 3
$ws.Range("D3").Value = "2.991.30"
$ws.Range("E3").Value = "  +0.65%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "'560.89"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("E6").Value = "  +3.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.65%  "

# Row 9
$ws.Range("D9").Value = "2.979.19"
$ws.Range("E9").Value = "  +0.28%  "

# Row 10
$ws.Range("E10").Value = "  +1.90%  "

# Row 11
$ws.Range("D11").Value = "'5.15"
$ws.Range("E11").Value = "  +5.07%  "

# Row 12
$ws.Range("E12").Value = "  +1.96%  "

# Row 13
$ws.Range("E13").Value = "  +1.61%  "

# Row 14
$ws.Range("D14").Value = "'33.68"
$ws.Range("E14").Value = "  +1.39%  "

# Row 15
$ws.Range("E15").Value = "  +1.59%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'7.30"
$ws.Range("E16").Value = "  +6.62%  "

# Row 17
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.486.05"
$ws.Range("E17").Value = "  +0.71%  "

# Row 18
$ws.Range("D18").Value = "2.990.44"
$ws.Range("E18").Value = "  +0.67%  "

# Row 19
$ws.Range("D19").Value = "59.306.25"
$ws.Range("E19").Value = "  +2.30%  "

# Row 20
$ws.Range("D20").Value = "'429.71"
$ws.Range("E20").Value = "  +1.48%  "

# Row 21
$ws.Range("D21").Value = "'13.59"
$ws.Range("E21").Value = "  +2.04%  "

# Row 22
$ws.Range("D22").Value = "'0.719"
$ws.Range("E22").Value = "  +4.09%  "

# Row 23
$ws.Range("E23").Value = "  +2.97%  "

# Row 24
$ws.Range("D24").Value = "'7.12"
$ws.Range("E24").Value = "  +1.12%  "

# Row 25
$ws.Range("D25").Value = "'80.31"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("E26").Value = "  -0.15%  "

# Row 27
$ws.Range("E27").Value = "  +9.37%  "

# Row 28
$ws.Range("E28").Value = "  +0.24%  "

# Row 29
$ws.Range("D29").Value = "'2.54"
$ws.Range("E29").Value = "  +1.10%  "

# Row 30
$ws.Range("D30").Value = "'7.85"
$ws.Range("E30").Value = "  +2.30%  "

# Row 31
$ws.Range("D31").Value = "'25.70"
$ws.Range("E31").Value = "  +0.74%  "

# Row 32
$ws.Range("D32").Value = "'6.10"
$ws.Range("E32").Value = "  -0.81%  "

# Row 33
$ws.Range("E33").Value = "  +1.50%  "

# Row 34
$ws.Range("E34").Value = "  +5.96%  "

# Row 35
$ws.Range("D35").Value = "'5.95"
$ws.Range("E35").Value = "  +4.36%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0759"
$ws.Range("E36").Value = "  +8.17%  "

# Row 37
$ws.Range("D37").Value = "'2.10"
$ws.Range("E37").Value = "  -1.47%  "

# Row 38
$ws.Range("D38").Value = "'48.82"
$ws.Range("E38").Value = "  +0.13%  "

# Row 39
$ws.Range("D39").Value = "'8.67"
$ws.Range("E39").Value = "  -1.19%  "

# Row 40
$ws.Range("D40").Value = "'2.74"
$ws.Range("E40").Value = "  +4.96%  "

# Row 41
$ws.Range("D41").Value = "'406.50"
$ws.Range("E41").Value = "  +6.88%  "

# Row 42
$ws.Range("D42").Value = "'0.0353"
$ws.Range("E42").Value = "  +0.06%  "

# Row 43
$ws.Range("D43").Value = "2.772.09"
$ws.Range("E43").Value = "  +2.17%  "

# Row 44
$ws.Range("E44").Value = "  -1.63%  "

# Row 45
$ws.Range("E45").Value = "  +3.67%  "

# Row 46
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").Value = "'123.04"
$ws.Range("E47").Value = "  +0.24%  "

# Row 48
$ws.Range("D48").Value = "'34.39"
$ws.Range("E48").Value = "  +19.21%  "

# Row 49
$ws.Range("D49").Value = "'0.110"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("D51").Value = "'23.48"
$ws.Range("E51").Value = "  -0.93%  "
